$d = $word.ActiveDocument

$replacements = @(
    @("908÷4=", "543÷8="),
    @("625÷5=", "862÷6="),
    @("689÷2=", "437÷5="),
    @("434÷2=", "767÷6="),
    @("391÷7=", "329÷8="),
    @("481÷2=", "511÷2="),
    @("824÷7=", "753÷9="),
    @("150÷3=", "869÷6="),
    @("831÷6=", "204÷4="),
    @("562÷6=", "280÷3="),
    @("874÷4=", "813÷2="),
    @("591÷4=", "911÷2="),
    @("859÷7=", "707÷8="),
    @("980÷6=", "351÷6="),
    @("444÷9=", "535÷5="),
    @("503÷7=", "965÷6="),
    @("447÷7=", "527÷5="),
    @("651÷5=", "747÷4="),
    @("806÷2=", "934÷3="),
    @("490÷2=", "875÷4="),
    @("383÷7=", "554÷4="),
    @("412÷4=", "509÷2="),
    @("909÷6=", "686÷4="),
    @("138÷4=", "532÷8="),
    @("964÷9=", "407÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
